$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# -- Cell address -> new text value (re-scrape refresh: timestamps + a few metrics) --

# Values that are safe to assign directly: Excel keeps dates/units/symbols as text
# because they do not parse as a recognised number/date literal.
$plainUpdates = [ordered]@{
    "E2" = "2026-02-15 23:18:56"
    "E3" = "2026-02-15 23:18:59"
    "G3" = "208 cm"
    "I3" = "3.2 mm"
    "E4" = "2026-02-15 23:19:01"
    "E5" = "2026-02-15 23:19:04"
    "I5" = "9.7 mm"
    "O5" = "-4.1 °C"
    "E6" = "2026-02-15 23:19:06"
    "J6" = "1015.5 hPa"
    "E7" = "2026-02-15 23:19:09"
    "E8" = "2026-02-15 23:19:12"
    "E9" = "2026-02-15 23:19:15"
    "N9" = "6.2 °C 22:55 TU"
    "O9" = "10.5 °C"
    "E10" = "2026-02-15 23:19:17"
    "O10" = "7.2 °C"
    "E11" = "2026-02-15 23:19:20"
    "O11" = "6.7 °C"
    "E12" = "2026-02-15 23:19:22"
    "N12" = "7.4 °C 22:58 TU"
    "E13" = "2026-02-15 23:19:25"
    "J13" = "1015.7 hPa"
    "O13" = "6.1 °C"
    "E14" = "2026-02-15 23:19:28"
    "O14" = "10.9 °C"
    "E15" = "2026-02-15 23:19:30"
    "E16" = "2026-02-15 23:19:33"
    "E17" = "2026-02-15 23:19:36"
    "E18" = "2026-02-15 23:19:39"
    "E19" = "2026-02-15 23:19:41"
    "E20" = "2026-02-15 23:19:44"
    "E21" = "2026-02-15 23:19:47"
    "E22" = "2026-02-15 23:19:50"
    "N22" = "-6.5 °C 22:53 TU"
    "E23" = "2026-02-15 23:19:53"
    "I23" = "6.3 mm"
    "E24" = "2026-02-15 23:19:56"
    "E25" = "2026-02-15 23:19:58"
    "E26" = "2026-02-15 23:20:01"
    "E27" = "2026-02-15 23:20:04"
    "E28" = "2026-02-15 23:20:07"
    "J28" = "1015.8 hPa"
    "O28" = "6.4 °C"
    "E29" = "2026-02-15 23:20:10"
    "E30" = "2026-02-15 23:20:12"
    "E31" = "2026-02-15 23:20:15"
    "O31" = "10.3 °C"
    "E32" = "2026-02-15 23:20:18"
    "O32" = "4.1 °C"
    "E33" = "2026-02-15 23:20:21"
    "E34" = "2026-02-15 23:20:24"
    "O34" = "1.5 °C"
    "E35" = "2026-02-15 23:20:27"
    "J35" = "1019.5 hPa"
    "O35" = "4.4 °C"
    "E36" = "2026-02-15 23:20:30"
    "E37" = "2026-02-15 23:20:33"
    "E38" = "2026-02-15 23:20:35"
    "E39" = "2026-02-15 23:20:38"
    "E40" = "2026-02-15 23:20:41"
    "J40" = "1016.4 hPa"
    "O40" = "8.3 °C"
    "E41" = "2026-02-15 23:20:44"
    "O41" = "12.8 °C"
    "E42" = "2026-02-15 23:20:47"
    "E43" = "2026-02-15 23:20:49"
    "E44" = "2026-02-15 23:20:52"
    "I44" = "6.9 mm"
    "O44" = "-3.6 °C"
    "E45" = "2026-02-15 23:20:55"
    "I45" = "5.5 mm"
    "O45" = "1.3 °C"
    "E46" = "2026-02-15 23:20:58"
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Bare "NN%" strings would otherwise be auto-parsed by Excel into a percentage
# number (e.g. "55%" -> 0.55). Source column stores these as literal text, so
# force the Text number format before assigning, then restore the original
# "General" format (the underlying value stays plain text either way).
$percentUpdates = [ordered]@{
    "H9" = "55%"
    "H10" = "73%"
    "H11" = "49%"
    "H12" = "61%"
    "H15" = "55%"
    "H17" = "42%"
    "H18" = "74%"
    "H20" = "64%"
    "H21" = "42%"
    "H28" = "61%"
    "H30" = "59%"
    "H33" = "44%"
    "H36" = "52%"
    "H37" = "58%"
    "H38" = "70%"
    "H40" = "44%"
    "H42" = "61%"
    "H43" = "70%"
}
foreach ($addr in $percentUpdates.Keys) {
    $cell = $ws.Range($addr)
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $percentUpdates[$addr]
    $cell.NumberFormat = $originalFormat
}

$wb.Save()
